# The document has two BTEC logos (image1.jpg) in the headers and two
# Pearson logos (image2.png) in the footers. This edit swaps the naming:
#   BTec_Logo-Orange pictures:    image1.jpg -> image2.jpg
#   PearsonLogo.png pictures:     image2.png -> image1.png
# across every header/footer in every section.

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    for ($i = 1; $i -le $sec.Headers.Count; $i++) {
        $hdr = $sec.Headers.Item($i)
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                $ishp = $shapes.Item($j)
                if ($ishp.AlternativeText -eq "BTec_Logo-Orange") {
                    $ishp.Name = "image2.jpg"
                }
            }
        }
    }
    for ($i = 1; $i -le $sec.Footers.Count; $i++) {
        $ftr = $sec.Footers.Item($i)
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($j = 1; $j -le $shapes.Count; $j++) {
                $ishp = $shapes.Item($j)
                if ($ishp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $ishp.Name = "image1.png"
                }
            }
        }
    }
}
